$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: merge the two runs that spell out "...CONSENTIMIENTO..."
# (currently split as "...CONSENT" + bookmark(_GoBack) + "IMIENTO...")
# into a single run, and drop the _GoBack bookmark from this spot
# (it gets re-created later, around "fecha").
# -----------------------------------------------------------------
$old1 = "CONSENTIMIENTO OTORGADO, PODRÁ REALIZARSE A TRAVÉS DE LA PRESENTACIÓN DE LA SOLICITUD RESPECTIVA EN: "
$found1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2)
if (-not $found1) {
    throw "Edit 1: could not find the CONSENTIMIENTO text"
}

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# -----------------------------------------------------------------
# Edit 2: the lone "{vigencia}" placeholder (the one right after
# "CULIACÁN, SINALOA A ", NOT the "({vigencia})" one earlier in the
# document) becomes "{" + bookmarked("fecha") + "}". Anchor on the
# unique surrounding context so we land on the right occurrence
# regardless of document order.
# -----------------------------------------------------------------
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("A {vigencia}.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundAnchor) {
    throw "Edit 2: could not find the CULIACÁN ... {vigencia}. context"
}

$target = $d.Range($anchor.End - 11, $anchor.End - 1)
if ($target.Text -ne "{vigencia}") {
    throw "Edit 2: unexpected target text [$($target.Text)]"
}

# Isolate the whole "{vigencia}" run from its neighbours using a
# scratch bookmark (bookmark boundaries stop Word from re-merging
# runs that happen to share identical formatting).
$scratchName = "ZZscratchZZ"
if ($d.Bookmarks.Exists($scratchName)) {
    $d.Bookmarks.Item($scratchName).Delete()
}
$null = $d.Bookmarks.Add($scratchName, $target)

$scratchRange = $d.Bookmarks.Item($scratchName).Range
$inner = $d.Range($scratchRange.Start + 1, $scratchRange.End - 1)
$inner.Text = "fecha"

$fechaBlock = $d.Bookmarks.Item($scratchName).Range
$goBackRange = $d.Range($fechaBlock.Start + 1, $fechaBlock.End - 1)
$d.Bookmarks.Item($scratchName).Delete()
$null = $d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "done"
